$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 355 to make room for a new weekly entry
# (this shifts the former rows 355-362 down to 357-364)
$ws.Rows("355:356").Insert()

# Row 355: new weekly entry - Primera
$ws.Cells.Item(355, 1).Value = 3
$ws.Cells.Item(355, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(355, 3).Value = "Coquimbo"
$ws.Cells.Item(355, 4).Value = 44448
$ws.Cells.Item(355, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(355, 5).Value = 5
$ws.Cells.Item(355, 6).Value = 100114014
$ws.Cells.Item(355, 7).Value = "Betarraga"
$ws.Cells.Item(355, 8).Value = "Sin especificar"
$ws.Cells.Item(355, 9).Value = "Primera"
$ws.Cells.Item(355, 10).Value = 3200
$ws.Cells.Item(355, 11).Value = 500
$ws.Cells.Item(355, 12).Value = 550
$ws.Cells.Item(355, 13).Value = 525
$ws.Cells.Item(355, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(355, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(355, 16).Value = 131
$ws.Cells.Item(355, 17).Value = 4
$ws.Cells.Item(355, 18).Value = "Hortaliza"

# Row 356: new weekly entry - Segunda
$ws.Cells.Item(356, 1).Value = 3
$ws.Cells.Item(356, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(356, 3).Value = "Coquimbo"
$ws.Cells.Item(356, 4).Value = 44448
$ws.Cells.Item(356, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(356, 5).Value = 5
$ws.Cells.Item(356, 6).Value = 100114014
$ws.Cells.Item(356, 7).Value = "Betarraga"
$ws.Cells.Item(356, 8).Value = "Sin especificar"
$ws.Cells.Item(356, 9).Value = "Segunda"
$ws.Cells.Item(356, 10).Value = 800
$ws.Cells.Item(356, 11).Value = 400
$ws.Cells.Item(356, 12).Value = 400
$ws.Cells.Item(356, 13).Value = 400
$ws.Cells.Item(356, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(356, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(356, 16).Value = 100
$ws.Cells.Item(356, 17).Value = 4
$ws.Cells.Item(356, 18).Value = "Hortaliza"
